$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "system" column F, duplicating the OPERACION header (col E) into F1,
# including its fill/style, by copying the cell.
$ws.Range("E1").Copy($ws.Range("F1"))

# Give the new column a sensible width (stored width ends up at 17).
$ws.Columns.Item(6).ColumnWidth = 16.14

# Re-apply the AutoFilter so its range grows from A1:E1 to A1:F1.
# (toggling the existing AutoFilter off then on again refreshes the range)
$ws.Range("A1:F1").AutoFilter()
$ws.Range("A1:F1").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$F`$1"
    }
}

# Move the active selection to F2, like a user would after adding the column.
$ws.Range("F2").Select()
